$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new rows are written as text (matching the existing sheet convention
# where every cell, including numeric-looking values, is stored as a string).
$dataRange = $ws.Range("A6:E13")
$dataRange.NumberFormat = "@"

$ws.Range("A6").Value = '2673867'
$ws.Range("B6").Value = 'PKT Mayo'
$ws.Range("C6").Value = '1'
$ws.Range("D6").Value = '$37.89'
$ws.Range("E6").Value = '$37.89'

$ws.Range("A7").Value = '3125531'
$ws.Range("B7").Value = 'PKT Sugar - (Raw)'
$ws.Range("C7").Value = '1'
$ws.Range("D7").Value = '$34.15'
$ws.Range("E7").Value = '$34.15'

$ws.Range("A8").Value = '7143223'
$ws.Range("B8").Value = 'Sausage - Vegan Patty'
$ws.Range("C8").Value = '1'
$ws.Range("D8").Value = '$94.93'
$ws.Range("E8").Value = '$94.93'

$ws.Range("A9").Value = '7468531'
$ws.Range("B9").Value = 'Gatorade Cool Blue'
$ws.Range("C9").Value = '1'
$ws.Range("D9").Value = '$29.40'
$ws.Range("E9").Value = '$29.40'

$ws.Range("A10").Value = '7468515'
$ws.Range("B10").Value = 'Gatorade Fruit Punch'
$ws.Range("C10").Value = '1'
$ws.Range("D10").Value = '$29.40'
$ws.Range("E10").Value = '$29.40'

$ws.Range("A11").Value = '2240263'
$ws.Range("B11").Value = 'Chobani - Strawberry'
$ws.Range("C11").Value = '1'
$ws.Range("D11").Value = '$15.27'
$ws.Range("E11").Value = '$15.27'

$ws.Range("A12").Value = '0543585'
$ws.Range("B12").Value = 'Chobani - Black Cherry'
$ws.Range("C12").Value = '1'
$ws.Range("D12").Value = '$13.25'
$ws.Range("E12").Value = '$13.25'

$ws.Range("A13").Value = '5756060'
$ws.Range("B13").Value = 'SABRA - Hummus Roasted Garlic With Pretzels'
$ws.Range("C13").Value = '1'
$ws.Range("D13").Value = '$28.77'
$ws.Range("E13").Value = '$28.77'

